$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 34 (pushes old rows 34-39 down by one, i.e. old empty rows and rows 38/39)
$ws.Rows.Item(34).Insert()

Write-Output "done"
